$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the trailing rows (69-77) that are being dropped from the table
$ws.Range("A69:C77").ClearContents() | Out-Null

# Rewrite the full data table (rows 2-68) with the updated BU mapping data
$ws.Cells.Item(2, 1).Value = 'Alan Bruno'
$ws.Cells.Item(2, 2).Value = 'alan.bruno@arzion.com'
$ws.Cells.Item(2, 3).Value = 'Club'
$ws.Cells.Item(3, 1).Value = 'Alongkron Rodthong'
$ws.Cells.Item(3, 2).Value = 'alongkron_ro@minor.com'
$ws.Cells.Item(3, 3).Value = 'IT'
$ws.Cells.Item(4, 1).Value = 'Anan Hayicheteh'
$ws.Cells.Item(4, 2).Value = 'anan_ha@minor.com'
$ws.Cells.Item(4, 3).Value = 'IT'
$ws.Cells.Item(5, 1).Value = 'Apiwat'
$ws.Cells.Item(5, 2).Value = 'apiwat.s@codemonday.com'
$ws.Cells.Item(5, 3).Value = 'Club'
$ws.Cells.Item(6, 1).Value = 'Ariel Orrino'
$ws.Cells.Item(6, 2).Value = 'ariel.orrino@arzion.com'
$ws.Cells.Item(6, 3).Value = 'Club'
$ws.Cells.Item(7, 1).Value = 'Deloitte'
$ws.Cells.Item(7, 2).Value = 'asukkaew@deloitte.com'
$ws.Cells.Item(7, 3).Value = 'Club'
$ws.Cells.Item(8, 1).Value = 'Atiwit'
$ws.Cells.Item(8, 2).Value = 'atiwit_wi@minor.com'
$ws.Cells.Item(8, 3).Value = 'IT'
$ws.Cells.Item(9, 1).Value = 'Boripat Kestin'
$ws.Cells.Item(9, 2).Value = 'boripat_ke@minor.com'
$ws.Cells.Item(9, 3).Value = 'IT'
$ws.Cells.Item(10, 1).Value = 'Chayanon Lormanometee'
$ws.Cells.Item(10, 2).Value = 'chayanon_lo@minor.com'
$ws.Cells.Item(10, 3).Value = 'Commercial'
$ws.Cells.Item(11, 1).Value = 'Deloitte'
$ws.Cells.Item(11, 2).Value = 'chchongchalearmpaibo@deloitte.com'
$ws.Cells.Item(11, 3).Value = 'Club'
$ws.Cells.Item(12, 1).Value = 'Chinnawat Phutthatham'
$ws.Cells.Item(12, 2).Value = 'chinnawat_ph@minor.com'
$ws.Cells.Item(12, 3).Value = 'IT'
$ws.Cells.Item(13, 1).Value = 'Deloitte'
$ws.Cells.Item(13, 2).Value = 'darchong@deloitte.com'
$ws.Cells.Item(13, 3).Value = 'Club'
$ws.Cells.Item(14, 1).Value = 'Arzion'
$ws.Cells.Item(14, 2).Value = 'dario.pereyra@arzion.com'
$ws.Cells.Item(14, 3).Value = 'Club'
$ws.Cells.Item(15, 1).Value = 'Dayin Promkotwong'
$ws.Cells.Item(15, 2).Value = 'dayin@hypcode.co'
$ws.Cells.Item(15, 3).Value = 'CA'
$ws.Cells.Item(16, 1).Value = 'Echo ITALOT'
$ws.Cells.Item(16, 2).Value = 'echo@italots.com'
$ws.Cells.Item(16, 3).Value = 'Club'
$ws.Cells.Item(17, 1).Value = 'Deloitte'
$ws.Cells.Item(17, 2).Value = 'echoyan@deloitte.com.cn'
$ws.Cells.Item(17, 3).Value = 'Club'
$ws.Cells.Item(18, 1).Value = 'German Orlando'
$ws.Cells.Item(18, 2).Value = 'german.orlando@arzion.com'
$ws.Cells.Item(18, 3).Value = 'Club'
$ws.Cells.Item(19, 1).Value = 'Deloitte'
$ws.Cells.Item(19, 2).Value = 'gteo@deloitte.com'
$ws.Cells.Item(19, 3).Value = 'Club'
$ws.Cells.Item(20, 1).Value = 'Guido Traversaro'
$ws.Cells.Item(20, 2).Value = 'guido.traversaro@arzion.com'
$ws.Cells.Item(20, 3).Value = 'Club'
$ws.Cells.Item(21, 1).Value = 'Deloitte'
$ws.Cells.Item(21, 2).Value = 'jduangjaidee@deloitte.com'
$ws.Cells.Item(21, 3).Value = 'Club'
$ws.Cells.Item(22, 1).Value = 'Jill Ji'
$ws.Cells.Item(22, 2).Value = 'jiji@minor.com'
$ws.Cells.Item(22, 3).Value = 'IT'
$ws.Cells.Item(23, 1).Value = 'Kamonchanok S.'
$ws.Cells.Item(23, 2).Value = 'kamonchanok_si@minor.com'
$ws.Cells.Item(23, 3).Value = 'FS'
$ws.Cells.Item(24, 1).Value = 'Khemmanij Tansui'
$ws.Cells.Item(24, 2).Value = 'khemmanij_ta@minor.com'
$ws.Cells.Item(24, 3).Value = 'IT'
$ws.Cells.Item(25, 1).Value = 'Khongsak Kawdettikhun'
$ws.Cells.Item(25, 2).Value = 'khongsak@hypcode.co'
$ws.Cells.Item(25, 3).Value = 'CA'
$ws.Cells.Item(26, 1).Value = 'Kitti Tongpraduppet'
$ws.Cells.Item(26, 2).Value = 'kitti.tongpraduppet@gmail.com'
$ws.Cells.Item(26, 3).Value = 'IT'
$ws.Cells.Item(27, 1).Value = 'Kittipong Balang'
$ws.Cells.Item(27, 2).Value = 'kittipong_ba@minor.com'
$ws.Cells.Item(27, 3).Value = 'IT'
$ws.Cells.Item(28, 1).Value = 'Deloitte'
$ws.Cells.Item(28, 2).Value = 'kkositanont@deloitte.com'
$ws.Cells.Item(28, 3).Value = 'Club'
$ws.Cells.Item(29, 1).Value = 'Kritsana Uttamang'
$ws.Cells.Item(29, 2).Value = 'kritsana@hypcode.co'
$ws.Cells.Item(29, 3).Value = 'CA'
$ws.Cells.Item(30, 1).Value = 'Lamai - Contract Admin'
$ws.Cells.Item(30, 2).Value = 'lamai_nu@minor.com'
$ws.Cells.Item(30, 3).Value = 'CA'
$ws.Cells.Item(31, 1).Value = 'Deloitte'
$ws.Cells.Item(31, 2).Value = 'mingsun@deloitte.com.cn'
$ws.Cells.Item(31, 3).Value = 'Club'
$ws.Cells.Item(32, 1).Value = 'Micheal Ye'
$ws.Cells.Item(32, 2).Value = 'mye@anantaraclub.com'
$ws.Cells.Item(32, 3).Value = 'IT'
$ws.Cells.Item(33, 1).Value = 'Nagorn - ITALOT'
$ws.Cells.Item(33, 2).Value = 'nagorn@italots.com'
$ws.Cells.Item(33, 3).Value = 'Club'
$ws.Cells.Item(34, 1).Value = 'Nattaphat Petprom'
$ws.Cells.Item(34, 2).Value = 'nattaphat_pe@minor.com'
$ws.Cells.Item(34, 3).Value = 'IT'
$ws.Cells.Item(35, 1).Value = 'Deloitte'
$ws.Cells.Item(35, 2).Value = 'nikang@deloitte.com.cn'
$ws.Cells.Item(35, 3).Value = 'Club'
$ws.Cells.Item(36, 1).Value = 'Deloitte'
$ws.Cells.Item(36, 2).Value = 'nyodkaew@deloitte.com'
$ws.Cells.Item(36, 3).Value = 'Club'
$ws.Cells.Item(37, 1).Value = 'Panachai'
$ws.Cells.Item(37, 2).Value = 'panachai_ch@anantaraclub.com'
$ws.Cells.Item(37, 3).Value = 'FS'
$ws.Cells.Item(38, 1).Value = 'Pasawish Imjumroon'
$ws.Cells.Item(38, 2).Value = 'pasawish@italots.com'
$ws.Cells.Item(38, 3).Value = 'IT'
$ws.Cells.Item(39, 1).Value = 'Oat - Marcom'
$ws.Cells.Item(39, 2).Value = 'pat_pa@minor.com'
$ws.Cells.Item(39, 3).Value = 'MarCom'
$ws.Cells.Item(40, 1).Value = 'Pearploy '
$ws.Cells.Item(40, 2).Value = 'pearploy_th@minor.com'
$ws.Cells.Item(40, 3).Value = 'Commercial'
$ws.Cells.Item(41, 1).Value = 'Pichalak Owchariyapitak'
$ws.Cells.Item(41, 2).Value = 'pichalak_ow@minor.com'
$ws.Cells.Item(41, 3).Value = 'IT'
$ws.Cells.Item(42, 1).Value = 'PIRAPOP THONGSANDEE'
$ws.Cells.Item(42, 2).Value = 'pirapop_th@anantaraclub.com'
$ws.Cells.Item(42, 3).Value = 'FS'
$ws.Cells.Item(43, 1).Value = 'Deloitte'
$ws.Cells.Item(43, 2).Value = 'pkeelawat@deloitte.com'
$ws.Cells.Item(43, 3).Value = 'Club'
$ws.Cells.Item(44, 1).Value = 'Deloitte'
$ws.Cells.Item(44, 2).Value = 'plertphati@deloitte.com'
$ws.Cells.Item(44, 3).Value = 'Club'
$ws.Cells.Item(45, 1).Value = 'Deloitte'
$ws.Cells.Item(45, 2).Value = 'psaejeam@deloitte.com'
$ws.Cells.Item(45, 3).Value = 'Club'
$ws.Cells.Item(46, 1).Value = 'Deloitte'
$ws.Cells.Item(46, 2).Value = 'pyordming@deloitte.com'
$ws.Cells.Item(46, 3).Value = 'Club'
$ws.Cells.Item(47, 1).Value = 'Deloitte'
$ws.Cells.Item(47, 2).Value = 'rbasheerahamed@deloitte.com'
$ws.Cells.Item(47, 3).Value = 'Club'
$ws.Cells.Item(48, 1).Value = 'Deloitte'
$ws.Cells.Item(48, 2).Value = 'rbundlukarn@deloitte.com'
$ws.Cells.Item(48, 3).Value = 'Club'
$ws.Cells.Item(49, 1).Value = 'Deloitte'
$ws.Cells.Item(49, 2).Value = 'rogallagher@deloitte.com'
$ws.Cells.Item(49, 3).Value = 'Club'
$ws.Cells.Item(50, 1).Value = 'Sahaschai - ITALOT'
$ws.Cells.Item(50, 2).Value = 'sahaschai@italots.com'
$ws.Cells.Item(50, 3).Value = 'Club'
$ws.Cells.Item(51, 1).Value = 'Sattaya Thomwan'
$ws.Cells.Item(51, 2).Value = 'sattaya_th@minor.com'
$ws.Cells.Item(51, 3).Value = 'IT'
$ws.Cells.Item(52, 1).Value = 'Deloitte'
$ws.Cells.Item(52, 2).Value = 'shawcxiao@deloitte.com.cn'
$ws.Cells.Item(52, 3).Value = 'Club'
$ws.Cells.Item(53, 1).Value = 'Win [CODEMONDAY]'
$ws.Cells.Item(53, 2).Value = 'sittiporn.k@codemonday.com'
$ws.Cells.Item(53, 3).Value = 'Club'
$ws.Cells.Item(54, 1).Value = 'Deloitte'
$ws.Cells.Item(54, 2).Value = 'sjanklan@deloitte.com'
$ws.Cells.Item(54, 3).Value = 'Club'
$ws.Cells.Item(55, 1).Value = 'Bobby Leong'
$ws.Cells.Item(55, 2).Value = 'sleong@anantaraclub.com'
$ws.Cells.Item(55, 3).Value = 'Club'
$ws.Cells.Item(56, 1).Value = 'Deloitte'
$ws.Cells.Item(56, 2).Value = 'slikitphatham@deloitte.com'
$ws.Cells.Item(56, 3).Value = 'Club'
$ws.Cells.Item(57, 1).Value = 'Andrew Ngan'
$ws.Cells.Item(57, 2).Value = 'sngan@anantaraclub.com'
$ws.Cells.Item(57, 3).Value = 'Club'
$ws.Cells.Item(58, 1).Value = 'Suchakree_si'
$ws.Cells.Item(58, 2).Value = 'suchakree_si@minor.com'
$ws.Cells.Item(58, 3).Value = 'IT'
$ws.Cells.Item(59, 1).Value = 'Surawut Issarolarn'
$ws.Cells.Item(59, 2).Value = 'surawut_is@minor.com'
$ws.Cells.Item(59, 3).Value = 'IT'
$ws.Cells.Item(60, 1).Value = 'Thungong C'
$ws.Cells.Item(60, 2).Value = 'thungong_ch@minor.com'
$ws.Cells.Item(60, 3).Value = 'IT'
$ws.Cells.Item(61, 1).Value = 'Tichanon Sankham'
$ws.Cells.Item(61, 2).Value = 'tichanon_sa@minor.com'
$ws.Cells.Item(61, 3).Value = 'IT'
$ws.Cells.Item(62, 1).Value = 'Tom√°s Bourgeois'
$ws.Cells.Item(62, 2).Value = 'tomas.bourgeois@arzion.com'
$ws.Cells.Item(62, 3).Value = 'Club'
$ws.Cells.Item(63, 1).Value = 'Deloitte'
$ws.Cells.Item(63, 2).Value = 'tvivitvorn@deloitte.com'
$ws.Cells.Item(63, 3).Value = 'Club'
$ws.Cells.Item(64, 1).Value = 'Deloitte'
$ws.Cells.Item(64, 2).Value = 'vping@deloitte.com.cn'
$ws.Cells.Item(64, 3).Value = 'Club'
$ws.Cells.Item(65, 1).Value = 'Whachiravit Thanyathanachot'
$ws.Cells.Item(65, 2).Value = 'whachiravit_th@minor.com'
$ws.Cells.Item(65, 3).Value = 'CA'
$ws.Cells.Item(66, 1).Value = 'Wissaroot Samart'
$ws.Cells.Item(66, 2).Value = 'wissaroot_sa@minor.com'
$ws.Cells.Item(66, 3).Value = 'IT'
$ws.Cells.Item(67, 1).Value = 'Oscar'
$ws.Cells.Item(67, 2).Value = 'wlee@anantaraclub.com'
$ws.Cells.Item(67, 3).Value = 'MarCom'
$ws.Cells.Item(68, 1).Value = 'Jasmine'
$ws.Cells.Item(68, 2).Value = 'yyeo@anantaraclub.com'
$ws.Cells.Item(68, 3).Value = 'Club'
